$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# 1. Update the model name in A8 ("FIRECLASS 240-2" -> "FC240-2")
$ws.Cells.Item(8, 1).Value2 = "FC240-2"

# 2. Populate the new row 9 values first (order matters so that number
#    formats such as quotePrefix carried by columns K/L/N/O are preserved
#    when we paste formats from row 8 afterwards instead of being reset).
$ws.Cells.Item(9, 1).Value2 = "FC718D"
$ws.Cells.Item(9, 2).Value2 = "Node1"
$ws.Cells.Item(9, 4).Value2 = "FIM"
$ws.Cells.Item(9, 5).Value2 = 16
$ws.Cells.Item(9, 6).Value2 = 0.276
$ws.Cells.Item(9, 7).Value2 = 0.426
$ws.Cells.Item(9, 8).Value2 = "410DIM"
$ws.Cells.Item(9, 9).Value2 = "Ancillary Conventional"
$ws.Cells.Item(9, 10).Value2 = "410DIM - 1"
$ws.Cells.Item(9, 11).Value2 = 0.277
$ws.Cells.Item(9, 12).Value2 = 0.431
$ws.Cells.Item(9, 13).Value2 = "Class B - 1 Spur"
$ws.Cells.Item(9, 14).Value2 = 0.277
$ws.Cells.Item(9, 15).Value2 = 0.431
$ws.Cells.Item(9, 17).Value2 = 0.26
$ws.Cells.Item(9, 18).Value2 = 0.46

# 3. Copy the formatting of row 8 down to row 9 (split around the blank
#    column P so no stray cell gets created there).
$ws.Range("A8:O8").Copy()
$ws.Range("A9:O9").PasteSpecial(-4122)
$ws.Range("Q8:R8").Copy()
$ws.Range("Q9:R9").PasteSpecial(-4122)

# 4. Update the selected cell on the sheet.
$ws.Activate()
$ws.Range("N8").Select()

Write-Host "done"
